# "Agrego test calculo huella"
# The "Gasoil" consumption row label is shortened from the verbose
# "Combustible consumido - Gasoil " to just "Gasoil", and the active
# selection ends up on that cell (B5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Gasoil"
$ws.Range("B5").Select()
